# Apply updated weight-estimation results to the Weights.xlsx workbook.
# All cells in this workbook hold static (pre-computed) values -- there are
# no formulas to recalculate, so each changed cell is simply overwritten
# with its new value, matching the target OOXML diff exactly.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# GLOBAL RESULTS sheet
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GLOBAL RESULTS")

$ws.Range("C2").Value  = 3200.0

$ws.Range("C6").Value  = 66515.68921150395
$ws.Range("C7").Value  = 65495.68921150396
$ws.Range("C8").Value  = 64520.21853515883

$ws.Range("C11").Value = 18683.04663529311

$ws.Range("C13").Value = 47832.642576210856
$ws.Range("C14").Value = 46812.642576210856
$ws.Range("C15").Value = 31512.642576210856
$ws.Range("C16").Value = 30567.947742152846
$ws.Range("C17").Value = 28505.32618821086

$ws.Range("C19").Value = 332.57844605800744

$ws.Range("C21").Value = 652296.0836059949
$ws.Range("C22").Value = 642293.3006059951
$ws.Range("C23").Value = 632727.2010978151

$ws.Range("C27").Value = 469077.98431999807
$ws.Range("C28").Value = 459075.201319998
$ws.Range("C29").Value = 309033.45631999813
$ws.Range("C30").Value = 299769.16472558316
$ws.Range("C31").Value = 279541.75706361793

$ws.Range("C33").Value = 3261.4804180347573

# ----------------------------------------------------------------------
# FUSELAGE sheet
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("FUSELAGE")

$ws.Range("C7").Value  = 6601.0
$ws.Range("D7").Value  = -13.182576672996126

$ws.Range("C8").Value  = 6600.0
$ws.Range("D8").Value  = -13.195728835293808

$ws.Range("C9").Value  = 7516.0
$ws.Range("D9").Value  = -1.1483481706164018

$ws.Range("C12").Value = 7296.0
$ws.Range("D12").Value = -4.041823876106608

# ----------------------------------------------------------------------
# WING sheet
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WING")

$ws.Range("C8").Value  = 6659.0
$ws.Range("D8").Value  = 31.370373110404227

$ws.Range("C14").Value = 5791.5714285714275
$ws.Range("D14").Value = 14.257531080790697

# ----------------------------------------------------------------------
# HORIZONTAL TAIL sheet
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")

$ws.Range("C8").Value  = 300.0
$ws.Range("D8").Value  = -60.54351310695172

$ws.Range("C9").Value  = 533.0
$ws.Range("D9").Value  = -29.898974953350898

$ws.Range("C10").Value = 491.66666666666663
$ws.Range("D10").Value = -35.3352020363931

# ----------------------------------------------------------------------
# LANDING GEARS sheet
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LANDING GEARS")

$ws.Range("C5").Value = 1863.0
$ws.Range("D5").Value = -8.115706147813832

$ws.Range("C6").Value = 2661.0
$ws.Range("D6").Value = 31.24213952800182

$ws.Range("C7").Value = 3066.0
$ws.Range("D7").Value = 51.21698601760751

$ws.Range("C8").Value = 2695.0
$ws.Range("D8").Value = 32.919040220956376

$ws.Range("C9").Value = 2571.25
$ws.Range("D9").Value = 26.815614904687934

Write-Host "Done applying Weights.xlsx updates."
